$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.882.04"
$ws.Range("E2").Value = "  +5.89%  "
# Row 3
$ws.Range("D3").Value = "2.450.20"
$ws.Range("E3").Value = "  +3.40%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").Value = "'575.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.38%  "
# Row 6
$ws.Range("D6").Value = "'145.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.91%  "
# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
# Row 8
$ws.Range("E8").Value = "  +1.70%  "
# Row 9
$ws.Range("D9").Value = "2.448.47"
$ws.Range("E9").Value = "  +3.38%  "
# Row 10
$ws.Range("E10").Value = "  +6.97%  "
# Row 11
$ws.Range("E11").Value = "  +1.25%  "
# Row 12
$ws.Range("E12").Value = "  +4.84%  "
# Row 13
$ws.Range("E13").Value = "  +5.61%  "
# Row 14
$ws.Range("D14").Value = "'26.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.84%  "
# Row 15
$ws.Range("E15").Value = "  +10.20%  "
# Row 16
$ws.Range("D16").Value = "2.877.85"
$ws.Range("E16").Value = "  +2.22%  "
# Row 17
$ws.Range("D17").Value = "62.386.71"
$ws.Range("E17").Value = "  +4.96%  "
# Row 18
$ws.Range("D18").Value = "2.444.19"
$ws.Range("E18").Value = "  +3.43%  "
# Row 19
$ws.Range("E19").Value = "  -3.76%  "
# Row 20
$ws.Range("D20").Value = "'10.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.93%  "
# Row 21
$ws.Range("D21").Value = "'326.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.36%  "
# Row 22
$ws.Range("E22").Value = "  +3.77%  "
# Row 23
$ws.Range("D23").Value = "'2.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +16.03%  "
# Row 24
$ws.Range("E24").Value = "  +0.04%  "
# Row 25
$ws.Range("E25").Value = "  +2.76%  "
# Row 26
$ws.Range("D26").Value = "'617.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.80%  "
# Row 27
$ws.Range("D27").Value = "'8.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.76%  "
# Row 28
$ws.Range("D28").Value = "0.0₃0989"
$ws.Range("E28").Value = "  +10.99%  "
# Row 29
$ws.Range("D29").Value = "2.531.17"
$ws.Range("E29").Value = "  +1.08%  "
# Row 30
$ws.Range("D30").Value = "'0.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "
# Row 31
$ws.Range("D31").Value = "'8.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.43%  "
# Row 32
$ws.Range("E32").Value = "  +10.61%  "
# Row 33
$ws.Range("D33").Value = "'0.139"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.67%  "
# Row 34
$ws.Range("E34").Value = "  +4.16%  "
# Row 35
$ws.Range("E35").Value = "  +6.58%  "
# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
# Row 37
$ws.Range("D37").Value = "'4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.96%  "
# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.373"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "
# Row 39
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'152.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.30%  "
# Row 40
$ws.Range("D40").Value = "'5.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.34%  "
# Row 41
$ws.Range("D41").Value = "'18.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.45%  "
# Row 42
$ws.Range("D42").Value = "'2.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +22.76%  "
# Row 43
$ws.Range("E43").Value = "  +8.97%  "
# Row 44
$ws.Range("D44").Value = "'42.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "
# Row 45
$ws.Range("E45").Value = "  +0.02%  "
# Row 46
$ws.Range("E46").Value = "  +0.48%  "
# Row 47
$ws.Range("D47").Value = "'144.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
# Row 48
$ws.Range("E48").Value = "  +4.03%  "
# Row 49
$ws.Range("D49").Value = "'20.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.22%  "
# Row 50
$ws.Range("D50").Value = "'0.601"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.52%  "
# Row 51
$ws.Range("D51").Value = "'0.0515"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
